# Update "想去人数" (F) / "最低票价" (G) figures that changed between scrapes.
# Same data is duplicated on the "展览" sheet (1) and the "全部类型" sheet (4).
$wb = $excel.ActiveWorkbook

$targetSheets = @(
    $wb.Worksheets.Item(1),
    $wb.Worksheets.Item(4)
)

foreach ($ws in $targetSheets) {
    $ws.Range("F2").Value = 1576
    $ws.Range("G2").Value = 70

    $ws.Range("F4").Value = 1031

    $ws.Range("F7").Value = 2675

    $ws.Range("F9").Value = 1701

    $ws.Range("F11").Value = 72

    $ws.Range("F12").Value = 572

    $ws.Range("F15").Value = 78
}

Write-Output "done"
